# Updated notes on income tax (slide 18 - "Joint Tax Assessment for Idris and his wife")
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(18)

# --- Table 1 ("Content Placeholder 3" / Aggregate -> Taxable Income) ---
$tbl1 = $s.Shapes.Item(1).Table
$tbl1.Cell(5, 2).Shape.TextFrame.TextRange.Text = "4,000"        # Self relief: 9,000 -> 4,000
$tbl1.Cell(11, 4).Shape.TextFrame.TextRange.Text = "(37,600)"    # Total Reliefs: (42,600) -> (37,600)
$tbl1.Cell(12, 4).Shape.TextFrame.TextRange.Text = "46,600"      # Taxable Income: 41,600 -> 46,600

# --- Table 2 ("Table 2" / Using Tax Rate -> Payable Tax) ---
$tbl2 = $s.Shapes.Item(2).Table
$tbl2.Cell(3, 1).Shape.TextFrame.TextRange.Text = "On the next 11,600 x 8%"  # was "On the next 6,600 x 8%"
$tbl2.Cell(3, 2).Shape.TextFrame.TextRange.Text = "928"                      # was 528
$tbl2.Cell(4, 4).Shape.TextFrame.TextRange.Text = "1,528"                    # Total Tax: 1,128 -> 1,528
$tbl2.Cell(8, 4).Shape.TextFrame.TextRange.Text = "RM 1,028"                 # Payable Tax: RM 628 -> RM 1,028

# --- TextBox 4 (explanatory note under the tables) ---
$note = $s.Shapes.Item(3).TextFrame.TextRange
$lastPara = $note.Paragraphs(3, 1)
$lastPara.Text = "Balance = 46,600 " + [char]0x2013 + " 35,000 = 11,600 @ 8% "
